$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers in I1 / J1, matching the formatting of the existing header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I2:J17
$data = @(
    @(3,6),
    @(6,6),
    @(6,8),
    @(1,6),
    @(1,5),
    @(1,6),
    @(1,6),
    @(1,5),
    @(1,3),
    @(1,5),
    @(1,5),
    @(1,7),
    @(1,5),
    @(6,7),
    @(5,6),
    @(6,7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
